$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("compounds")
$ws2 = $wb.Worksheets.Item("biomarkers")

# --- "compounds" sheet -----------------------------------------------
# Row 2 = Open Targets Platform: source_version "2024.06" -> "2024.09"
# "2024.09" parses as a number, so Value would store it numerically.
# Build it with TEXT() and paste-special as values to force plain text,
# matching the original (unformatted) string cell.
$c = $ws1.Range("E2")
$c.Formula = '=TEXT(2024.09,"0.00")'
$c.Copy()
$c.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Row 3 = NCI Thesaurus: source_version "24.07e" -> "24.09e"
$ws1.Range("E3").Value = "24.09e"

# --- "biomarkers" sheet ------------------------------------------------
# Row 3 = Mitelman Database: source_version "v20240715" -> "v20241015"
$ws2.Range("E3").Value = "v20241015"

# Row 2 (CIViC) source_version cell (E2) stays blank / skipped - no change.

# The active tab moves from "compounds" to "biomarkers", selection E3
$ws2.Activate()
$ws2.Range("E3").Select()
